$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("county-year")
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 3072
$ws.Range("V5").Value = 100
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 3072
$ws.Range("V6").Value = 100
$ws.Range("B7").Value = 2034
$ws.Range("C7").Value = 66.2109375
$ws.Range("S7").Value = 3056
$ws.Range("T7").Value = 99.479164123535156
$ws.Range("U7").Value = 16
$ws.Range("V7").Value = 0.52083331346511841
$ws.Range("B8").Value = 2020
$ws.Range("C8").Value = 65.755210876464844
$ws.Range("S8").Value = 3027
$ws.Range("T8").Value = 98.53515625
$ws.Range("U8").Value = 45
$ws.Range("V8").Value = 1.46484375
$ws.Range("S9").Value = 3072
$ws.Range("T9").Value = 100
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0

$ws = $wb.Worksheets.Item("county-overall")
$ws.Range("A2").Value = 4054
$ws.Range("B2").Value = 16.495767593383789
$ws.Range("R2").Value = 9155
$ws.Range("S2").Value = 37.251789093017578
$ws.Range("T2").Value = 15421
$ws.Range("U2").Value = 62.748210906982422

$ws = $wb.Worksheets.Item("point-year")
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1362620
$ws.Range("V5").Value = 100
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 1362620
$ws.Range("V6").Value = 100
$ws.Range("B7").Value = 915810
$ws.Range("C7").Value = 67.209495544433594
$ws.Range("S7").Value = 1356702
$ws.Range("T7").Value = 99.565689086914063
$ws.Range("U7").Value = 5918
$ws.Range("V7").Value = 0.43431037664413452
$ws.Range("B8").Value = 904424
$ws.Range("C8").Value = 66.373893737792969
$ws.Range("S8").Value = 1335431
$ws.Range("T8").Value = 98.004653930664063
$ws.Range("U8").Value = 27189
$ws.Range("V8").Value = 1.9953471422195435
$ws.Range("S9").Value = 1362620
$ws.Range("T9").Value = 100
$ws.Range("U9").Value = 0
$ws.Range("V9").Value = 0

$ws = $wb.Worksheets.Item("point-overall")
$ws.Range("A2").Value = 1820234
$ws.Range("B2").Value = 16.69792366027832
$ws.Range("R2").Value = 4054753
$ws.Range("S2").Value = 37.196292877197266
$ws.Range("T2").Value = 6846207
$ws.Range("U2").Value = 62.803707122802734
